# Update Data by bot, scripted by HH
# Applies the 688560 balance-sheet row-2 data refresh:
#  - DATE_TYPE_CODE (J2) and REPORT_DATE (N2) are text fields; use the
#    leading apostrophe so Excel keeps them as text instead of coercing
#    the numeric-looking / date-looking literal into a number/date.
#  - The remaining touched columns are plain numeric measures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "'001"
$ws.Range("N2").Value = "'2018-12-31 00:00:00"

$ws.Range("O2").Value = 985777199.4
$ws.Range("P2").Value = 143627512.87
$ws.Range("Q2").Value = 112909979.79
$ws.Range("R2").Value = -3.7078082524
$ws.Range("S2").Value = 350226580.97
$ws.Range("T2").Value = 42.1040956159
$ws.Range("U2").Value = 102225813.13
$ws.Range("V2").Value = 1.3213439175
$ws.Range("W2").Value = 443593298.87
$ws.Range("X2").Value = 241143906.98
$ws.Range("Y2").Value = 15.5558326774

$ws.Range("AB2").Value = 542183900.53
$ws.Range("AC2").Value = 19.5048442296
$ws.Range("AD2").Value = 19.2850140234
$ws.Range("AE2").Value = 19.0174214632
$ws.Range("AF2").Value = 208.6449631147
$ws.Range("AG2").Value = 44.9993466211
